$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect what it actually documents now.
$ws.Name = "campaign_details"

# The old sanity-test table (rows 1-5, columns A:C - Username/Password/
# Result with sample login rows) is being trimmed down to a single
# "owner" column. Drop the stale hyperlink collection before the
# underlying cells disappear.
$ws.Hyperlinks.Delete()

# Rows 3-5 (the extra sample rows) go away entirely...
$ws.Rows("3:5").Delete()

# ...and so do columns B:C (Password/Result), leaving just column A.
$ws.Columns("B:C").Delete()

# Strip the border/fill formatting that decorated the old table so the
# two remaining cells fall back to plain styling.
$ws.Cells.ClearFormats()

# New, smaller data set: a header ("owner") plus the single remaining
# campaign-owner e-mail address.
$ws.Range("A1").Value = "owner"
$ws.Range("A2").Value = "abhone@convirza.com"

# Re-create the mailto hyperlink on the owner e-mail cell.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:abhone@convirza.com")

# Column A now only needs to fit the e-mail address.
$ws.Columns("A").AutoFit()

# Leave the selection where the editor last left it.
$ws.Range("E13").Select() | Out-Null

Write-Host "campaigns_page.xlsx trimmed to campaign_details owner column"
